$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the auto-update strings: clear the content of the four rows
# (A20:B23) that held the "new version available" / "update app" texts.
# The cell formatting (fill/alignment/border) stays untouched.
$ws.Range("A20:B23").ClearContents()

# B23 previously carried a style (xf 6) that only differed from the
# style used by the rest of the column (xf 5) by a redundant
# "applyFont" flag pointing at the default font. Re-apply the same
# alignment so the cell collapses onto the shared style, matching its
# neighbours now that its distinguishing text is gone.
$ws.Range("B23").HorizontalAlignment = -4108

# Row 21 no longer needs the taller 30pt height that accommodated the
# two-line "new version available" prompt - let it size back to the
# sheet's default row height.
$ws.Rows(21).AutoFit()

# Update the surviving selection so it no longer points at the removed
# row 21 text; it now rests on row 17.
$ws.Range("A17").Select()
